$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1870109324866007
$ws.Range("D2").Value = 0.1531424273940765
$ws.Range("E2").Value = 0.1318313254792209
$ws.Range("F2").Value = 1.425286788478957
$ws.Range("G2").Value = 0.8679767137632979
$ws.Range("H2").Value = 0.8949953293570729
$ws.Range("I2").Value = 0.9939850508759704
$ws.Range("J2").Value = 0.1480452987242324
$ws.Range("L2").Value = 0.1750774811426949
$ws.Range("M2").Value = 1.158243764853978
$ws.Range("N2").Value = 1.86477819059678
$ws.Range("O2").Value = 3.554511214537541
$ws.Range("C3").Value = 0.1881661551333451
$ws.Range("D3").Value = 0.1541696090682905
$ws.Range("E3").Value = 0.1333230058917891
$ws.Range("F3").Value = 1.41739764474832
$ws.Range("G3").Value = 0.8550740402857571
$ws.Range("H3").Value = 0.893492746379124
$ws.Range("I3").Value = 0.9878313419564435
$ws.Range("J3").Value = 0.1500106438187074
$ws.Range("L3").Value = 0.1771055688953211
$ws.Range("M3").Value = 1.057171418215162
$ws.Range("N3").Value = 1.714121077930599
$ws.Range("O3").Value = 3.52364790160297
$ws.Range("C4").Value = 0.1889535385620746
$ws.Range("D4").Value = 0.1548534623564137
$ws.Range("E4").Value = 0.134292022877387
$ws.Range("F4").Value = 1.413346828502725
$ws.Range("G4").Value = 0.8477037080632925
$ws.Range("H4").Value = 0.8929961380174092
$ws.Range("I4").Value = 0.9846129173704838
$ws.Range("J4").Value = 0.1512820724479695
$ws.Range("L4").Value = 0.1784171030564607
$ws.Range("M4").Value = 0.9949182052387187
$ws.Range("N4").Value = 1.621788000599253
$ws.Range("O4").Value = 3.506745862013446
$ws.Range("C5").Value = 0.189294044498368
$ws.Range("D5").Value = 0.1551455172503573
$ws.Range("E5").Value = 0.1347002588391364
$ws.Range("F5").Value = 1.411895502172499
$ws.Range("G5").Value = 0.844838911367674
$ws.Range("H5").Value = 0.8929009400609687
$ws.Range("I5").Value = 0.983442164179273
$ws.Range("J5").Value = 0.1518164482523285
$ws.Range("L5").Value = 0.1789682296887358
$ws.Range("M5").Value = 0.9695034479248079
$ws.Range("N5").Value = 1.584209141978903
$ws.Range("O5").Value = 3.500373021441817
$ws.Range("C6").Value = 0.1893517717682727
$ws.Range("D6").Value = 0.1551948211536178
$ws.Range("E6").Value = 0.1347688525336874
$ws.Range("F6").Value = 1.411666553835786
$ws.Range("G6").Value = 0.8443715860657193
$ws.Range("H6").Value = 0.8928916072529631
$ws.Range("I6").Value = 0.9832562634802002
$ws.Range("J6").Value = 0.1519061625966744
$ws.Range("L6").Value = 0.1790607503966029
$ws.Range("M6").Value = 0.9652806473714861
$ws.Range("N6").Value = 1.577972214032371
$ws.Range("O6").Value = 3.499345909244994
$ws.Range("C7").Value = 0.1889580512067823
$ws.Range("D7").Value = 0.1548573469231478
$ws.Range("E7").Value = 0.1342974744266523
$ws.Range("F7").Value = 1.413326448030446
$ws.Range("G7").Value = 0.8476645110927024
$ws.Range("H7").Value = 0.892994420128403
$ws.Range("I7").Value = 0.9845965582085796
$ws.Range("J7").Value = 0.151289213412459
$ws.Range("L7").Value = 0.1784244682650833
$ws.Range("M7").Value = 0.9945756357271591
$ws.Range("N7").Value = 1.621280999999442
$ws.Range("O7").Value = 3.506657831142263
$ws.Range("C8").Value = 0.1873930547038185
$ws.Range("D8").Value = 0.1534855734162299
$ws.Range("E8").Value = 0.1323346300177124
$ws.Range("F8").Value = 1.42240196904929
$ws.Range("G8").Value = 0.8634132327259607
$ws.Range("H8").Value = 0.8943888571529328
$ws.Range("I8").Value = 0.9917470059485396
$ws.Range("J8").Value = 0.148709515266058
$ws.Range("L8").Value = 0.1757630090859141
$ws.Range("M8").Value = 1.123436016048672
$ws.Range("N8").Value = 1.812798836154315
$ws.Range("O8").Value = 3.543444339655196
$ws.Range("C9").Value = 0.1849432210457209
$ws.Range("D9").Value = 0.1512168128050781
$ws.Range("E9").Value = 0.1289070797579956
$ws.Range("F9").Value = 1.446496598370516
$ws.Range("G9").Value = 0.8986843391666355
$ws.Range("H9").Value = 0.9005020818585194
$ws.Range("I9").Value = 1.010215665430245
$ws.Range("J9").Value = 0.1441644861513127
$ws.Range("L9").Value = 0.1710697530507344
$ws.Range("M9").Value = 1.374477020366953
$ws.Range("N9").Value = 2.189539052771863
$ws.Range("O9").Value = 3.631849307835807
$ws.Range("C10").Value = 0.1835203621297552
$ws.Range("D10").Value = 0.1498060889532198
$ws.Range("E10").Value = 0.1266459299115432
$ws.Range("F10").Value = 1.468048220775216
$ws.Range("G10").Value = 0.9272886778289546
$ws.Range("H10").Value = 0.9070534082744359
$ws.Range("I10").Value = 1.026503284392163
$ws.Range("J10").Value = 0.1411388693486426
$ws.Range("L10").Value = 0.1679420114195578
$ws.Range("M10").Value = 1.557776552854776
$ws.Range("N10").Value = 2.466830846553194
$ws.Range("O10").Value = 3.706751261024976
$ws.Range("C11").Value = 0.1829548583697047
$ws.Range("D11").Value = 0.1492197867103755
$ws.Range("E11").Value = 0.1256731061978025
$ws.Range("F11").Value = 1.478690901926612
$ws.Range("G11").Value = 0.9408893979228594
$ws.Range("H11").Value = 0.9104812456105265
$ws.Range("I11").Value = 1.034505177373589
$ws.Range("J11").Value = 0.1398306234713875
$ws.Range("L11").Value = 0.166588649299718
$ws.Range("M11").Value = 1.640889894638008
$ws.Range("N11").Value = 2.593043748012406
$ws.Range("O11").Value = 3.742994980887374
$ws.Range("C12").Value = 0.1827524673649776
$ws.Range("D12").Value = 0.1490057309044239
$ws.Range("E12").Value = 0.1253127487340375
$ws.Range("F12").Value = 1.482841732706092
$ws.Range("G12").Value = 0.9461244614826398
$ws.Range("H12").Value = 0.9118436234183775
$ws.Range("I12").Value = 1.037620589216971
$ws.Range("J12").Value = 0.1393450320989278
$ws.Range("L12").Value = 0.1660861559690678
$ws.Range("M12").Value = 1.67232124933426
$ws.Range("N12").Value = 2.640843434668795
$ws.Range("O12").Value = 3.757032043104005
$ws.Range("C13").Value = 0.182795533239279
$ws.Range("D13").Value = 0.1490514775805991
$ws.Range("E13").Value = 0.1253900009081796
$ws.Range("F13").Value = 1.481942408737581
$ws.Range("G13").Value = 0.9449932245744321
$ws.Range("H13").Value = 0.911547350431448
$ws.Range("I13").Value = 1.036945836280395
$ws.Range("J13").Value = 0.1394491761833603
$ws.Range("L13").Value = 0.1661939322349655
$ws.Range("M13").Value = 1.665553853607562
$ws.Range("N13").Value = 2.630548741760322
$ws.Range("O13").Value = 3.753995014774091
$ws.Range("C14").Value = 0.1829379720164468
$ws.Range("D14").Value = 0.1492020166454608
$ws.Range("E14").Value = 0.1256432983981322
$ws.Range("F14").Value = 1.479029974307096
$ws.Range("G14").Value = 0.941318390083751
$ws.Range("H14").Value = 0.9105920401790399
$ws.Range("I14").Value = 1.034759774943751
$ws.Range("J14").Value = 0.1397904768124361
$ws.Range("L14").Value = 0.16654710850655
$ws.Range("M14").Value = 1.643476627026274
$ws.Range("N14").Value = 2.5969761713759
$ws.Range("O14").Value = 3.744143557047096
$ws.Range("C15").Value = 0.1830267503400549
$ws.Range("D15").Value = 0.1492952630814095
$ws.Range("E15").Value = 0.125799496374241
$ws.Range("F15").Value = 1.477261742040824
$ws.Range("G15").Value = 0.9390784911985577
$ws.Range("H15").Value = 0.9100152617232595
$ws.Range("I15").Value = 1.033431855440654
$ws.Range("J15").Value = 0.1400008117454328
$ws.Range("L15").Value = 0.1667647410873858
$ws.Range("M15").Value = 1.629948148742244
$ws.Range("N15").Value = 2.576412593985935
$ws.Range("O15").Value = 3.738149939432049
$ws.Range("C16").Value = 0.183558964024634
$ws.Range("D16").Value = 0.1498455201121836
$ws.Range("E16").Value = 0.1267106294426439
$ws.Range("F16").Value = 1.467369580003108
$ws.Range("G16").Value = 0.9264116926480597
$ws.Range("H16").Value = 0.9068383950211398
$ws.Range("I16").Value = 1.02599227086634
$ws.Range("J16").Value = 0.1412257385439639
$ws.Range("L16").Value = 0.1680318551732993
$ws.Range("M16").Value = 1.552339203106257
$ws.Range("N16").Value = 2.458583573774774
$ws.Range("O16").Value = 3.704426344556907
$ws.Range("C17").Value = 0.1839063980359512
$ws.Range("D17").Value = 0.150197279219018
$ws.Range("E17").Value = 0.1272838735491639
$ws.Range("F17").Value = 1.461515929431158
$ws.Range("G17").Value = 0.9187918455401416
$ws.Range("H17").Value = 0.9050041078444764
$ws.Range("I17").Value = 1.021580139866586
$ws.Range("J17").Value = 0.1419946513248924
$ws.Range("L17").Value = 0.1688269846306598
$ws.Range("M17").Value = 1.504657225685364
$ws.Range("N17").Value = 2.386314247384519
$ws.Range("O17").Value = 3.684294089173648
$ws.Range("C18").Value = 0.1841139289646776
$ws.Range("D18").Value = 0.1504048208308397
$ws.Range("E18").Value = 0.1276188386376642
$ws.Range("F18").Value = 1.458228006546307
$ws.Range("G18").Value = 0.9144644912829989
$ws.Range("H18").Value = 0.9039912072403808
$ws.Range("I18").Value = 1.019098171218815
$ws.Range("J18").Value = 0.1424433205954163
$ws.Range("L18").Value = 0.1692908598576155
$ws.Range("M18").Value = 1.477206525317172
$ws.Range("N18").Value = 2.344753845975902
$ws.Range("O18").Value = 3.672918826407454
$ws.Range("C19").Value = 0.1841855173086486
$ws.Range("D19").Value = 0.1504759874967618
$ws.Range("E19").Value = 0.1277331535518869
$ws.Range("F19").Value = 1.457128328639484
$ws.Range("G19").Value = 0.9130088306436903
$ws.Range("H19").Value = 0.9036554939357302
$ws.Range("I19").Value = 1.018267396496483
$ws.Range("J19").Value = 0.1425963328749174
$ws.Range("L19").Value = 0.1694490427064501
$ws.Range("M19").Value = 1.467907937107
$ws.Range("N19").Value = 2.330683554456868
$ws.Range("O19").Value = 3.669102435376885
$ws.Range("C20").Value = 0.1838686166090397
$ws.Range("D20").Value = 0.1501592937685885
$ws.Range("E20").Value = 0.1272223072972984
$ws.Range("F20").Value = 1.462130889814574
$ws.Range("G20").Value = 0.9195972585628454
$ws.Range("H20").Value = 0.9051950105040021
$ws.Range("I20").Value = 1.02204404609035
$ws.Range("J20").Value = 0.1419121356119417
$ws.Range("L20").Value = 0.168741665056241
$ws.Range("M20").Value = 1.509735692306322
$ws.Range("N20").Value = 2.394006749746097
$ws.Range("O20").Value = 3.686416058329655
$ws.Range("C21").Value = 0.1828958153565594
$ws.Range("D21").Value = 0.1491575835968675
$ws.Range("E21").Value = 0.1255686807961176
$ws.Range("F21").Value = 1.479882151176255
$ws.Range("G21").Value = 0.9423954763701659
$ws.Range("H21").Value = 0.9108708923313316
$ws.Range("I21").Value = 1.035399559528244
$ws.Range("J21").Value = 0.1396899620875165
$ws.Range("L21").Value = 0.1664431007436766
$ws.Range("M21").Value = 1.649962403327791
$ws.Range("N21").Value = 2.60683713414403
$ws.Range("O21").Value = 3.747028689895558
$ws.Range("C22").Value = 0.1823285369645973
$ws.Range("D22").Value = 0.1485493251896735
$ws.Range("E22").Value = 0.1245347476107961
$ws.Range("F22").Value = 1.492187064205197
$ws.Range("G22").Value = 0.957789577897131
$ws.Range("H22").Value = 0.9149553588649724
$ws.Range("I22").Value = 1.044625196789369
$ws.Range("J22").Value = 0.1382948485771756
$ws.Range("L22").Value = 0.1649991209356045
$ws.Range("M22").Value = 1.741363287305575
$ws.Range("N22").Value = 2.74596460130283
$ws.Range("O22").Value = 3.788463288091577
$ws.Range("C23").Value = 0.182625037461591
$ws.Range("D23").Value = 0.1488697197933391
$ws.Range("E23").Value = 0.1250822916926353
$ws.Range("F23").Value = 1.485555313429671
$ws.Range("G23").Value = 0.9495281936882805
$ws.Range("H23").Value = 0.9127411057724828
$ws.Range("I23").Value = 1.039655802795423
$ws.Range("J23").Value = 0.1390342069690449
$ws.Range("L23").Value = 0.1657644666216074
$ws.Range("M23").Value = 1.692604341539052
$ws.Range("N23").Value = 2.671708394257905
$ws.Range("O23").Value = 3.766182184998627
$ws.Range("C24").Value = 0.1838856733380183
$ws.Range("D24").Value = 0.1501764504485443
$ws.Range("E24").Value = 0.1272501245815016
$ws.Range("F24").Value = 1.461852625071828
$ws.Range("G24").Value = 0.919232965015695
$ws.Range("H24").Value = 0.9051085736498123
$ws.Range("I24").Value = 1.021834143937156
$ws.Range("J24").Value = 0.1419494203732183
$ws.Range("L24").Value = 0.1687802170266073
$ws.Range("M24").Value = 1.507439834735891
$ws.Range("N24").Value = 2.390529005832832
$ws.Range("O24").Value = 3.685456096086341
$ws.Range("C25").Value = 0.1855397260978506
$ws.Range("D25").Value = 0.1517855429380539
$ws.Range("E25").Value = 0.1297891832212454
$ws.Range("F25").Value = 1.43930311706481
$ws.Range("G25").Value = 0.8886711919881378
$ws.Range("H25").Value = 0.8984865595687666
$ws.Range("I25").Value = 1.00474254211661
$ws.Range("J25").Value = 0.1453389934388647
$ws.Range("L25").Value = 0.1722831259490492
$ws.Range("M25").Value = 1.306754546491177
$ws.Range("N25").Value = 2.087517568983628
$ws.Range("O25").Value = 3.631849307835807
